$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Obrigatorio) from "N" to "S" for rows 2 through 8
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 5).Value = "S"
}
